# Update images to match GitHub default dark theme
#
# This script rewrites the raw OOXML of the document to apply a set of
# color / shadow tweaks to the FIRST logo drawing (the dark "Decorator"
# graphic group) while leaving the second logo drawing (the light
# "bg1" rectangle graphic group) untouched.
#
# Changes applied only inside the first <mc:AlternateContent> ... 
# </mc:AlternateContent> block (the "AC1" group, which contains both the
# mc:Choice (DrawingML) and mc:Fallback (VML) representations of the
# first logo):
#   1. wp:anchor wp14:editId 6CF1326C -> 57E624B6
#   2. Rectangle solid fill color 22272E -> 0D1117 (DrawingML + VML fallback)
#   3. "Decorator" text color A02B93 (themeColor accent5) -> D769CA (explicit, no theme link)
#   4. "Decorator" text outline accent5 lumMod/lumOff 20000/80000 -> lumMod 75000 only
#   5. "TEXT" text shadow accent1 lumMod/lumOff 60000/40000 removed (alpha kept)

$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML

$acStartTag = "<mc:AlternateContent>"
$acEndTag = "</mc:AlternateContent>"

$ac1Start = $xml.IndexOf($acStartTag)
$ac1EndTagStart = $xml.IndexOf($acEndTag)
$ac1End = $ac1EndTagStart + $acEndTag.Length

$before = $xml.Substring(0, $ac1Start)
$ac1 = $xml.Substring($ac1Start, $ac1End - $ac1Start)
$after = $xml.Substring($ac1End)

# 1. editId on the wp:anchor (unique within the document)
$ac1 = $ac1.Replace('wp14:editId="6CF1326C"', 'wp14:editId="57E624B6"')

# 2. Rectangle / v:rect fill color (dark background square)
$ac1 = $ac1.Replace('<a:srgbClr val="22272E"/>', '<a:srgbClr val="0D1117"/>')
$ac1 = $ac1.Replace('fillcolor="#22272e"', 'fillcolor="#0d1117"')

# 3. "Decorator" run/paragraph color: drop the theme link, use explicit new color
$ac1 = $ac1.Replace('<w:color w:val="A02B93" w:themeColor="accent5"/>', '<w:color w:val="D769CA"/>')

# 4. "Decorator" text outline: collapse lumMod/lumOff pair into a single lumMod
$ac1 = $ac1.Replace('<w14:schemeClr w14:val="accent5"><w14:lumMod w14:val="20000"/><w14:lumOff w14:val="80000"/></w14:schemeClr>', '<w14:schemeClr w14:val="accent5"><w14:lumMod w14:val="75000"/></w14:schemeClr>')

# 5. "TEXT" shadow: remove the lumMod/lumOff pair, keep alpha
$ac1 = $ac1.Replace('<w14:alpha w14:val="50000"/><w14:lumMod w14:val="60000"/><w14:lumOff w14:val="40000"/>', '<w14:alpha w14:val="50000"/>')

$newXml = $before + $ac1 + $after

$d.Content.InsertXML($newXml)
